$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep exact text formatting (no numeric coercion / trailing-zero loss)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "52.340.83"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "2.908.48"
$ws.Range("E3").Value = "  +3.73%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "352.93"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").Value = "112.30"
$ws.Range("E6").Value = "  +0.85%  "
$ws.Range("D7").Value = "0.559"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "0.633"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").Value = "39.92"
$ws.Range("E10").Value = "  -1.23%  "
$ws.Range("D11").Value = "0.0865"
$ws.Range("E11").Value = "  +3.01%  "
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").Value = "19.85"
$ws.Range("E13").Value = "  -1.00%  "
$ws.Range("D14").Value = "7.81"
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("D15").Value = "3.367.15"
$ws.Range("E15").Value = "  +3.60%  "
$ws.Range("D16").Value = "1.00"
$ws.Range("E16").Value = "  +6.19%  "
$ws.Range("D17").Value = "2.902.19"
$ws.Range("E17").Value = "  +3.24%  "
$ws.Range("D18").Value = "52.372.22"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("D19").Value = "7.61"
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("D20").Value = "3.31"
$ws.Range("E20").Value = "  +3.73%  "
$ws.Range("D21").Value = "14.18"
$ws.Range("E21").Value = "  +3.98%  "
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("D23").Value = "70.98"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("D24").Value = "269.94"
$ws.Range("E24").Value = "  +0.37%  "
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").Value = "0.172"
$ws.Range("E26").Value = "  +6.05%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "26.80"
$ws.Range("E27").Value = "  +2.32%  "
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("D29").Value = "10.66"
$ws.Range("E29").Value = "  +2.54%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "6.37"
$ws.Range("E30").Value = "  +12.67%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "6.65"
$ws.Range("E31").Value = "  +8.20%  "
$ws.Range("D32").Value = "37.92"
$ws.Range("E32").Value = "  -2.19%  "
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("D34").Value = "0.0982"
$ws.Range("E34").Value = "  +10.68%  "
$ws.Range("D35").Value = "53.37"
$ws.Range("E35").Value = "  +1.81%  "
$ws.Range("D36").Value = "0.0452"
$ws.Range("E36").Value = "  +1.89%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D39").Value = "18.87"
$ws.Range("E39").Value = "  +0.38%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "2.07"
$ws.Range("E40").Value = "  +2.97%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "2.86"
$ws.Range("E41").Value = "  +14.32%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "0.117"
$ws.Range("E42").Value = "  +1.29%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "23.52"
$ws.Range("E43").Value = "  +6.38%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "121.56"
$ws.Range("E44").Value = "  +0.92%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "2.62"
$ws.Range("E45").Value = "  +7.92%  "
$ws.Range("E46").Value = "  -0.58%  "
$ws.Range("D47").Value = "3.54"
$ws.Range("E47").Value = "  +3.44%  "
$ws.Range("D48").Value = "2.197.19"
$ws.Range("E48").Value = "  +3.98%  "
$ws.Range("D49").Value = "0.265"
$ws.Range("E49").Value = "  +22.27%  "
$ws.Range("D50").Value = "0.0338"
$ws.Range("E50").Value = "  +9.78%  "
$ws.Range("D51").Value = "0.968"
$ws.Range("E51").Value = "  +1.90%  "
